$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 10-15 previously only had the match metadata (index/home/away) filled
# in; this adds the xG_home, xG_away, goals_home and goals_away figures for
# those six matches, completing the table the same way it is already filled
# in for rows 2-9.
#
# The values being written look numeric but the rest of the sheet stores
# these figures as text (shared strings) rather than numbers, so the
# destination cells are pre-formatted as Text before the values are written
# and the formatting is reset back to Normal afterwards (leaving the cells
# with no explicit style, matching the rest of the sheet).
$ws.Range("D10:D15").NumberFormat = "@"
$ws.Range("E10:E15").NumberFormat = "@"
$ws.Range("F10:F15").NumberFormat = "@"
$ws.Range("G10:G15").NumberFormat = "@"

$ws.Range("D10").Value = "0.753035"
$ws.Range("E10").Value = "1.49985"
$ws.Range("F10").Value = "0"
$ws.Range("G10").Value = "3"

$ws.Range("D11").Value = "1.98287"
$ws.Range("E11").Value = "1.1214"
$ws.Range("F11").Value = "3"
$ws.Range("G11").Value = "1"

$ws.Range("D12").Value = "1.61081"
$ws.Range("E12").Value = "2.74594"
$ws.Range("F12").Value = "1"
$ws.Range("G12").Value = "3"

$ws.Range("D13").Value = "1.31797"
$ws.Range("E13").Value = "3.25581"
$ws.Range("F13").Value = "1"
$ws.Range("G13").Value = "0"

$ws.Range("D14").Value = "1.26331"
$ws.Range("E14").Value = "0.779401"
$ws.Range("F14").Value = "2"
$ws.Range("G14").Value = "1"

$ws.Range("D15").Value = "1.23629"
$ws.Range("E15").Value = "2.0312"
$ws.Range("F15").Value = "1"
$ws.Range("G15").Value = "2"

$ws.Range("D10:D15").Style = "Normal"
$ws.Range("E10:E15").Style = "Normal"
$ws.Range("F10:F15").Style = "Normal"
$ws.Range("G10:G15").Style = "Normal"
